$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '57.906.88'
$ws.Range('E2').Value = '  +2.51%  '

# Row 3
$ws.Range('D3').Value = '3.066.86'
$ws.Range('E3').Value = '  +1.59%  '

# Row 4
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '517.42'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.52%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.00'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.59%  '

# Row 7
$ws.Range('E7').Value = '  +0.04%  '

# Row 8
$ws.Range('E8').Value = '  +1.55%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.29'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.48%  '

# Row 10
$ws.Range('E10').Value = '  +0.08%  '

# Row 11
$ws.Range('E11').Value = '  +2.22%  '

# Row 12
$ws.Range('D12').Value = '3.591.24'
$ws.Range('E12').Value = '  +1.84%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.44'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.91%  '

# Row 15
$ws.Range('E15').Value = '  +0.98%  '

# Row 16
$ws.Range('D16').Value = '57.919.71'
$ws.Range('E16').Value = '  +2.48%  '

# Row 17
$ws.Range('D17').Value = '3.064.41'
$ws.Range('E17').Value = '  +1.66%  '

# Row 18
$ws.Range('E18').Value = '  +2.69%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.87'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.91%  '

# Row 20
$ws.Range('E20').Value = '  +0.88%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '332.47'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.01%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.09%  '

# Row 23
$ws.Range('E23').Value = '  +0.00%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.41'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.17%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.171'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.71%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.13%  '

# Row 27
$ws.Range('D27').Value = '0.0₃0905'
$ws.Range('E27').Value = '  -4.08%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.48'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.45%  '

# Row 29
$ws.Range('E29').Value = '  +6.32%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.83'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.82%  '

# Row 31
$ws.Range('E31').Value = '  +2.60%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.71'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.50%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '154.93'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.35%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.55'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.59%  '

# Row 35
$ws.Range('E35').Value = '  +3.01%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '26.99'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.31%  '

# Row 37
$ws.Range('E37').Value = '  +3.65%  '

# Row 38
$ws.Range('E38').Value = '  +2.23%  '

# Row 39
$ws.Range('D39').Value = '3.106.59'
$ws.Range('E39').Value = '  +1.63%  '

# Row 40
$ws.Range('E40').Value = '  +3.70%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.49'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.11%  '

# Row 42
$ws.Range('E42').Value = '  -0.05%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.656'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.06%  '

# Row 44
$ws.Range('D44').Value = '2.266.27'
$ws.Range('E44').Value = '  +2.72%  '

# Row 45
$ws.Range('E45').Value = '  +7.44%  '

# Row 46
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.38'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.05%  '

# Row 47
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '20.82'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.94%  '

# Row 48
$ws.Range('E48').Value = '  +1.37%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.93'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.48%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.743'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +9.79%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '256.12'
$ws.Range('D51').Style = 'Normal'
